$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Ror2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.156959333333335
$ws.Range("H2").Value = 27.470878
$ws.Range("I2").Value = 0.969469463764299
$ws.Range("J2").Value = 0.9694694637642989
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.028814
$ws.Range("N2").Value = 0.086442
$ws.Range("O2").Value = 0.009404897244300481
$ws.Range("P2").Value = 0.009404897244300481
$ws.Range("Q2").Value = 0.2638486262306667
$ws.Range("R2").Value = 2.374637636076
$ws.Range("S2").Value = 0.009117760688190321
$ws.Range("T2").Value = 0.00911776068819032
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Ror2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 9.156959333333335
$ws.Range("H3").Value = 27.470878
$ws.Range("I3").Value = 0.969469463764299
$ws.Range("J3").Value = 0.9694694637642989
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.294804333333333
$ws.Range("N3").Value = 6.884412999999999
$ws.Range("O3").Value = 0.7490247432073112
$ws.Range("P3").Value = 0.7490247432073114
$ws.Range("Q3").Value = 21.01342995829044
$ws.Range("R3").Value = 189.120869624614
$ws.Range("S3").Value = 0.7261566161433837
$ws.Range("T3").Value = 0.7261566161433838
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Ror2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 9.156959333333335
$ws.Range("H4").Value = 27.470878
$ws.Range("I4").Value = 0.969469463764299
$ws.Range("J4").Value = 0.9694694637642989
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7401046666666667
$ws.Range("N4").Value = 2.220314
$ws.Range("O4").Value = 0.2415703595483883
$ws.Range("P4").Value = 0.2415703595483883
$ws.Range("Q4").Value = 6.777108335076891
$ws.Range("R4").Value = 60.99397501569201
$ws.Range("S4").Value = 0.2341950869327249
$ws.Range("T4").Value = 0.2341950869327249
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Ror2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.288371
$ws.Range("H5").Value = 0.865113
$ws.Range("I5").Value = 0.03053053623570109
$ws.Range("J5").Value = 0.03053053623570109
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.028814
$ws.Range("N5").Value = 0.086442
$ws.Range("O5").Value = 0.009404897244300481
$ws.Range("P5").Value = 0.009404897244300481
$ws.Range("Q5").Value = 0.008309121994000001
$ws.Range("R5").Value = 0.074782097946
$ws.Range("S5").Value = 0.0002871365561101611
$ws.Range("T5").Value = 0.0002871365561101611
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Ror2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.288371
$ws.Range("H6").Value = 0.865113
$ws.Range("I6").Value = 0.03053053623570109
$ws.Range("J6").Value = 0.03053053623570109
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.294804333333333
$ws.Range("N6").Value = 6.884412999999999
$ws.Range("O6").Value = 0.7490247432073112
$ws.Range("P6").Value = 0.7490247432073114
$ws.Range("Q6").Value = 0.6617550204076665
$ws.Range("R6").Value = 5.955795183668999
$ws.Range("S6").Value = 0.02286812706392751
$ws.Range("T6").Value = 0.02286812706392752
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Ror2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.288371
$ws.Range("H7").Value = 0.865113
$ws.Range("I7").Value = 0.03053053623570109
$ws.Range("J7").Value = 0.03053053623570109
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7401046666666667
$ws.Range("N7").Value = 2.220314
$ws.Range("O7").Value = 0.2415703595483883
$ws.Range("P7").Value = 0.2415703595483883
$ws.Range("Q7").Value = 0.2134247228313333
$ws.Range("R7").Value = 1.920822505482
$ws.Range("S7").Value = 0.007375272615663408
$ws.Range("T7").Value = 0.007375272615663409
